# Adds a "Total" column (T) with row totals, a new "Outros" data row (row 7),
# and a new "Total" summary row (row 8) to the sp2014_c analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" column header
$ws.Range("T1").Value = "Total"

# Row totals for the existing category rows (2-6)
$ws.Range("T2").Value = 82592
$ws.Range("T3").Value = 8938
$ws.Range("T4").Value = 37898
$ws.Range("T5").Value = 13332
$ws.Range("T6").Value = 50912

# New category row: "Outros"
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 7336
$ws.Range("C7").Value = 348
$ws.Range("D7").Value = 555
$ws.Range("E7").Value = 2601
$ws.Range("F7").Value = 3184
$ws.Range("G7").Value = 3411
$ws.Range("H7").Value = 3618
$ws.Range("I7").Value = 3932
$ws.Range("J7").Value = 4302
$ws.Range("K7").Value = 4777
$ws.Range("L7").Value = 5276
$ws.Range("M7").Value = 5558
$ws.Range("N7").Value = 5508
$ws.Range("O7").Value = 5240
$ws.Range("P7").Value = 5463
$ws.Range("Q7").Value = 6433
$ws.Range("R7").Value = 19784
$ws.Range("S7").Value = 626
$ws.Range("T7").Value = 87952

# New summary row: "Total" (column totals)
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 8221
$ws.Range("C8").Value = 568
$ws.Range("D8").Value = 811
$ws.Range("E8").Value = 3192
$ws.Range("F8").Value = 3916
$ws.Range("G8").Value = 4504
$ws.Range("H8").Value = 5367
$ws.Range("I8").Value = 6558
$ws.Range("J8").Value = 8428
$ws.Range("K8").Value = 11729
$ws.Range("L8").Value = 15675
$ws.Range("M8").Value = 20472
$ws.Range("N8").Value = 23408
$ws.Range("O8").Value = 25489
$ws.Range("P8").Value = 27634
$ws.Range("Q8").Value = 31505
$ws.Range("R8").Value = 83342
$ws.Range("S8").Value = 805
$ws.Range("T8").Value = 281624
